# Loan RBI, Variable Instalments
# Inserts a new (blank) column before column N on the "Repayment Schedule"
# sheet - shifting the old N/O/P columns ("Late" / "Heading" / "Outstanding")
# one column to the right - and makes "Repayment Schedule" the active sheet
# with L15 selected (previously "Transactions" was the active sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")

# Insert a blank column before column N; existing N/O/P data shifts to O/P/Q.
$ws.Columns("N").Insert()

# Give the freshly inserted column the same width as its left neighbour (M).
$ws.Columns("N").ColumnWidth = $ws.Columns("M").ColumnWidth

# Make "Repayment Schedule" the active sheet/tab with L15 as the selection
# (the "Transactions" sheet was previously the active tab at I3).
$ws.Activate()
$ws.Range("L15").Select()
